# Generate Report for Handoff
# Marks the 4d46be1d-... and bc12ed2a-... files as handed off (zh-cn & de-de)
# and records that the handback file version is stale for both locales.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$overviewDate = "2016-08-30 14:32:08"
$zhHandoffDate = "2016-08-30 14:31:56"
$deHandoffDate = "2016-08-30 14:32:08"

function Get-ErrorDetail($fileName) {
    return "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c1bec57991b48259c657d4132b06fbeb55a25e9/e2e/$fileName.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d4b739942051c2618b44abcc89e056f22655320f/e2e/$fileName.md."
}

$file1 = "4d46be1d-59b0-46ce-b6ed-c013e59cfeb2"
$file2 = "bc12ed2a-3bd0-4803-9ccf-0fd133318b58"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 4 -> 4d46be1d file
$wsOverview.Range("E4").Value = $statusReady
$wsOverview.Range("F4").Value = $statusReady
$wsOverview.Range("G4").Value = $overviewDate

# Row 5 -> bc12ed2a file
$wsOverview.Range("E5").Value = $statusReady
$wsOverview.Range("F5").Value = $statusReady
$wsOverview.Range("G5").Value = $overviewDate

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 4 -> 4d46be1d file
$wsZh.Range("C4").Value = $statusReady
$wsZh.Range("H4").Value = $zhHandoffDate
$wsZh.Range("P4").Value = Get-ErrorDetail($file1)

# Row 5 -> bc12ed2a file
$wsZh.Range("C5").Value = $statusReady
$wsZh.Range("H5").Value = $zhHandoffDate
$wsZh.Range("P5").Value = Get-ErrorDetail($file2)

$wsZh.Range("P1").ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

# Row 4 -> 4d46be1d file
$wsDe.Range("C4").Value = $statusReady
$wsDe.Range("H4").Value = $deHandoffDate
$wsDe.Range("P4").Value = Get-ErrorDetail($file1)

# Row 5 -> bc12ed2a file
$wsDe.Range("C5").Value = $statusReady
$wsDe.Range("H5").Value = $deHandoffDate
$wsDe.Range("P5").Value = Get-ErrorDetail($file2)

$wsDe.Range("P1").ColumnWidth = 39.17
